$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRF = 13.97171428571429

for ($r = 29; $r -le 47; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
